$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Condition Groupers (2)")
$ws2 = $wb.Worksheets.Item("Condition Groupers (3)")

# --- Condition Groupers (3): append / refresh 3 data rows (new "Carbapenem-resistant
#     Organisms" entries) copied from the matching rows already present on the
#     "Condition Groupers (2)" sheet (rows 64:66 -> rows 2:4) -----------------------
$ws1.Range("A64:F66").Copy($ws2.Range("A2:F4"))

# Match the source row height on the newly created row 4 (rows 2/3 already carried
# the correct height because they previously held data).
$ws2.Rows("4").RowHeight = $ws1.Rows("66").RowHeight

# The editor's cursor ends up on D2 after keying in the new rows.
$ws2.Activate()
$ws2.Range("D2").Select()

# --- Condition Groupers (2): the editor rezoomed / rescrolled while reviewing the
#     source rows (64:66) that were copied above, and left that range selected -----
$ws1.Activate()
$w = $excel.ActiveWindow
$w.Zoom = 120
$ws1.Range("A64:XFD66").Select()
